$wb = $excel.ActiveWorkbook

$secComm   = $wb.Worksheets.Item("SEC_Comm")
$secProc   = $wb.Worksheets.Item("SEC_Processes")
$finalDmd  = $wb.Worksheets.Item("FINAL_DEMAD_PRC")
$demand    = $wb.Worksheets.Item("DEMAND")

# ---------------------------------------------------------------------
# SEC_Comm (sheet1): rows 7-9, new "ELC_GRID_RES"/"ELC_GRID" commodity
# set rows plus the "ELC_FIN" demand row gets its Cset/CommName filled.
# ---------------------------------------------------------------------

# Csets column (C7, C8) is filled in first, picking up style from the
# blank s=34 cells in the same rows (F7/F8).
$secComm.Range("F7").Copy()
$secComm.Range("C7").PasteSpecial(-4122)
$secComm.Range("C7").Value = "ELC_GRID_RES"

$secComm.Range("F8").Copy()
$secComm.Range("C8").PasteSpecial(-4122)
$secComm.Range("C8").Value = "ELC_GRID"

# CommName column (B7, B8) is filled in next; it loses its old style
# entirely (copied from a genuinely style-less cell elsewhere in the
# workbook).
$demand.Range("L8").Copy()
$secComm.Range("B7").PasteSpecial(-4122)
$secComm.Range("B7").Value = "NRG"

$demand.Range("L8").Copy()
$secComm.Range("B8").PasteSpecial(-4122)
$secComm.Range("B8").Value = "NRG"

# Row 9: styles are untouched, just fill in the Cset/CommName values.
$secComm.Range("B9").Value = "DEM"
$secComm.Range("C9").Value = "ELC_FIN"

# ---------------------------------------------------------------------
# SEC_Processes (sheet2): row 7 - new "DMD" process, Tact unit filled in.
# ---------------------------------------------------------------------

$secProc.Range("C7").Copy()
$secProc.Range("D7").PasteSpecial(-4122)
$secProc.Range("D7").Value = "ELC_FIN_DEM"

$demand.Range("L8").Copy()
$secProc.Range("B7").PasteSpecial(-4122)
$secProc.Range("B7").Value = "DMD"

# ---------------------------------------------------------------------
# FINAL_DEMAD_PRC (sheet3): row 7 gains Comm-IN/Comm-OUT entries, row 8
# keeps its Comm-OUT, and the old placeholder rows 9/10 disappear.
# ---------------------------------------------------------------------

$finalDmd.Range("D8").Copy()
$finalDmd.Range("D7").PasteSpecial(-4122)
$finalDmd.Range("E7").PasteSpecial(-4122)
$finalDmd.Range("D7").Value = "ELC_GRID_RES"
$finalDmd.Range("E7").Value = "ELC_FIN"

$finalDmd.Range("B7").Value = "ELC_FIN_DEM"

$finalDmd.Range("D8").Value = "ELC_GRID"

$finalDmd.Range("D9").ClearFormats()
$finalDmd.Range("D9").ClearContents()
$finalDmd.Range("E10").ClearFormats()
$finalDmd.Range("E10").ClearContents()

# ---------------------------------------------------------------------
# Sheet-view bookkeeping: each sheet's last active cell moved, and the
# DEMAND sheet was also scrolled/zoomed out. SEC_Comm must be
# re-selected last so it stays the active (tabSelected) sheet.
# ---------------------------------------------------------------------

$finalDmd.Select()
$finalDmd.Range("E7").Select()

$secProc.Select()
$secProc.Range("D13").Select()

$demand.Select()
$excel.ActiveWindow.Zoom = 83
$demand.Range("C8").Select()

$secComm.Select()
$secComm.Range("B10").Select()
